# Auto-generated edit script for betclever_predictions workbook update
$wb = $excel.ActiveWorkbook

# ---- Sheet: Home win ----
$ws = $wb.Worksheets.Item("Home win")
$ws.Cells.Item(2, 1).Value = "08-05-2025 22:00"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(2, 4).Value = "Bodo/Glimt - Tottenham"
$ws.Cells.Item(2, 5).Value = 90
$ws.Cells.Item(2, 6).Value = 2.9
$ws.Cells.Item(3, 1).Value = "08-05-2025 03:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(3, 4).Value = "Universidad De Chile - Estudiantes L.P."
$ws.Cells.Item(3, 5).Value = 80
$ws.Cells.Item(3, 6).Value = 2.3
$ws.Cells.Item(4, 1).Value = "08-05-2025 05:00"
$ws.Cells.Item(4, 2).Value = "WORLD"
$ws.Cells.Item(4, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(4, 4).Value = "Sporting Cristal - Bolívar"
$ws.Cells.Item(4, 5).Value = 70
$ws.Cells.Item(4, 6).Value = 2.05
$ws.Cells.Item(5, 1).Value = "08-05-2025 03:30"
$ws.Cells.Item(5, 2).Value = "WORLD"
$ws.Cells.Item(5, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(5, 4).Value = "Mushuc Runa SC - Cruzeiro"
$ws.Cells.Item(5, 5).Value = 70
$ws.Cells.Item(5, 6).Value = 2.2
$ws.Cells.Item(6, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(6, 2).Value = "ITALY"
$ws.Cells.Item(6, 3).Value = "SERIE B"
$ws.Cells.Item(6, 4).Value = "Juve Stabia - Reggiana"
$ws.Cells.Item(6, 5).Value = 73.3
$ws.Cells.Item(6, 6).Value = 2.25
$ws.Cells.Item(7, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(7, 2).Value = "ITALY"
$ws.Cells.Item(7, 3).Value = "SERIE B"
$ws.Cells.Item(7, 4).Value = "Sassuolo - Catanzaro"
$ws.Cells.Item(7, 5).Value = 70
$ws.Cells.Item(7, 6).Value = 1.75
$ws.Cells.Item(8, 1).Value = "09-05-2025 20:30"
$ws.Cells.Item(8, 2).Value = "FRANCE"
$ws.Cells.Item(8, 3).Value = "NATIONAL 1"
$ws.Cells.Item(8, 4).Value = "Rouen - Gobelins"
$ws.Cells.Item(8, 5).Value = 73.3
$ws.Cells.Item(8, 6).Value = 2.45
$ws.Cells.Item(9, 1).Value = "09-05-2025 21:00"
$ws.Cells.Item(9, 2).Value = "HUNGARY"
$ws.Cells.Item(9, 3).Value = "NB I"
$ws.Cells.Item(9, 4).Value = "Gyori ETO FC - Fehérvár FC"
$ws.Cells.Item(9, 5).Value = 73.3
$ws.Cells.Item(9, 6).Value = 1.83
$ws.Cells.Item(10, 1).Value = "09-05-2025 21:00"
$ws.Cells.Item(10, 2).Value = "WORLD"
$ws.Cells.Item(10, 3).Value = "AFRICA CUP OF NATIONS U20"
$ws.Cells.Item(10, 4).Value = "South Africa U20 - Zambia U20"
$ws.Cells.Item(10, 5).Value = 80
$ws.Cells.Item(10, 6).Value = 2.3

# ---- Sheet: Away win ----
$ws = $wb.Worksheets.Item("Away win")
$ws.Cells.Item(2, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(2, 2).Value = "LITHUANIA"
$ws.Cells.Item(2, 3).Value = "A LYGA"
$ws.Cells.Item(2, 4).Value = "Džiugas Telšiai - Kauno Žalgiris"
$ws.Cells.Item(2, 5).Value = 73.3
$ws.Cells.Item(2, 6).Value = 1.7

# ---- Sheet: Btts ----
$ws = $wb.Worksheets.Item("Btts")
$ws.Cells.Item(2, 1).Value = "08-05-2025 03:30"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(2, 4).Value = "Atletico Grau - Gremio"
$ws.Cells.Item(2, 5).Value = 88
$ws.Cells.Item(2, 6).Value = 2.1
$ws.Cells.Item(3, 1).Value = "09-05-2025 01:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(3, 4).Value = "Deportes Iquique - Atletico-MG"
$ws.Cells.Item(3, 5).Value = 88
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(4, 1).Value = "08-05-2025 20:00"
$ws.Cells.Item(4, 2).Value = "DENMARK"
$ws.Cells.Item(4, 3).Value = "DBU POKALEN"
$ws.Cells.Item(4, 4).Value = "FC Copenhagen - Viborg"
$ws.Cells.Item(4, 5).Value = 76.7
$ws.Cells.Item(4, 6).Value = 1.8
$ws.Cells.Item(5, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(5, 2).Value = "ITALY"
$ws.Cells.Item(5, 3).Value = "SERIE B"
$ws.Cells.Item(5, 4).Value = "Modena - Brescia"
$ws.Cells.Item(5, 5).Value = 83.3
$ws.Cells.Item(5, 6).Value = 1.91
$ws.Cells.Item(6, 1).Value = "09-05-2025 12:30"
$ws.Cells.Item(6, 2).Value = "AUSTRALIA"
$ws.Cells.Item(6, 3).Value = "VICTORIA NPL"
$ws.Cells.Item(6, 4).Value = "Heidelberg United - Port Melbourne"
$ws.Cells.Item(6, 5).Value = 83.3
$ws.Cells.Item(6, 6).Value = 1.7
$ws.Cells.Item(7, 1).Value = "09-05-2025 19:00"
$ws.Cells.Item(7, 2).Value = "POLAND"
$ws.Cells.Item(7, 3).Value = "EKSTRAKLASA"
$ws.Cells.Item(7, 4).Value = "Motor Lublin - Piast Gliwice"
$ws.Cells.Item(7, 5).Value = 80
$ws.Cells.Item(7, 6).Value = 1.77
$ws.Cells.Item(8, 1).Value = "09-05-2025 02:00"
$ws.Cells.Item(8, 2).Value = "USA"
$ws.Cells.Item(8, 3).Value = "USL CHAMPIONSHIP"
$ws.Cells.Item(8, 4).Value = "North Carolina - Orange County SC"
$ws.Cells.Item(8, 5).Value = 84
$ws.Cells.Item(8, 6).Value = 1.74
$ws.Range("A9:F11").ClearContents()

# ---- Sheet: Over_Under ----
$ws = $wb.Worksheets.Item("Over_Under")
$ws.Cells.Item(2, 1).Value = "08-05-2025 05:00"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(2, 4).Value = "Sporting Cristal - Bolívar"
$ws.Cells.Item(2, 5).Value = 86.7
$ws.Cells.Item(2, 6).Value = 1.62
$ws.Cells.Item(2, 7).Value = 60
$ws.Cells.Item(2, 8).Value = 2.38
$ws.Cells.Item(3, 1).Value = "09-05-2025 01:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(3, 4).Value = "Deportes Iquique - Atletico-MG"
$ws.Cells.Item(3, 5).Value = 80
$ws.Cells.Item(3, 6).Value = 1.95
$ws.Cells.Item(3, 7).Value = 53.3
$ws.Cells.Item(3, 8).Value = 3.4
$ws.Cells.Item(4, 1).Value = "09-05-2025 05:00"
$ws.Cells.Item(4, 2).Value = "WORLD"
$ws.Cells.Item(4, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(4, 4).Value = "Universitario - Independiente Del Valle"
$ws.Cells.Item(4, 5).Value = 85
$ws.Cells.Item(4, 6).Value = 2.5
$ws.Cells.Item(4, 7).Value = 75
$ws.Cells.Item(4, 8).Value = 5
$ws.Cells.Item(5, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(5, 2).Value = "AUSTRIA"
$ws.Cells.Item(5, 3).Value = "2. LIGA"
$ws.Cells.Item(5, 4).Value = "Admira Wacker - Voitsberg"
$ws.Cells.Item(5, 5).Value = 80
$ws.Cells.Item(5, 6).Value = 1.76
$ws.Cells.Item(5, 7).Value = 60
$ws.Cells.Item(5, 8).Value = 2.92
$ws.Cells.Item(6, 1).Value = "09-05-2025 19:31"
$ws.Cells.Item(6, 2).Value = "FINLAND"
$ws.Cells.Item(6, 3).Value = "VEIKKAUSLIIGA"
$ws.Cells.Item(6, 4).Value = "Ilves - Haka"
$ws.Cells.Item(6, 5).Value = 85
$ws.Cells.Item(6, 6).Value = 1.48
$ws.Cells.Item(6, 7).Value = 70
$ws.Cells.Item(6, 8).Value = 2.2
$ws.Cells.Item(7, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(7, 2).Value = "GERMANY"
$ws.Cells.Item(7, 3).Value = "REGIONALLIGA - NORD"
$ws.Cells.Item(7, 4).Value = "Werder Bremen II - SSV Jeddeloh"
$ws.Cells.Item(7, 5).Value = 80
$ws.Cells.Item(7, 6).Value = 1.48
$ws.Cells.Item(7, 7).Value = 65
$ws.Cells.Item(7, 8).Value = 2.2
$ws.Cells.Item(8, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(8, 2).Value = "GERMANY"
$ws.Cells.Item(8, 3).Value = "REGIONALLIGA - NORDOST"
$ws.Cells.Item(8, 4).Value = "Hertha BSC II - BFC Dynamo"
$ws.Cells.Item(8, 5).Value = 85
$ws.Cells.Item(8, 6).Value = 1.5
$ws.Cells.Item(8, 7).Value = 73.8
$ws.Cells.Item(8, 8).Value = 2.3

# ---- Sheet: EV Home win ----
$ws = $wb.Worksheets.Item("EV Home win")
$ws.Cells.Item(2, 1).Value = "08-05-2025 22:00"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(2, 4).Value = "Bodo/Glimt - Tottenham"
$ws.Cells.Item(2, 5).Value = 90
$ws.Cells.Item(2, 6).Value = 2.9
$ws.Cells.Item(2, 7).Value = 1.61
$ws.Cells.Item(3, 1).Value = "08-05-2025 03:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(3, 4).Value = "Universidad De Chile - Estudiantes L.P."
$ws.Cells.Item(3, 5).Value = 80
$ws.Cells.Item(3, 6).Value = 2.3
$ws.Cells.Item(3, 7).Value = 0.84
$ws.Cells.Item(4, 1).Value = "08-05-2025 05:00"
$ws.Cells.Item(4, 2).Value = "WORLD"
$ws.Cells.Item(4, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(4, 4).Value = "Sporting Cristal - Bolívar"
$ws.Cells.Item(4, 5).Value = 70
$ws.Cells.Item(4, 6).Value = 2.05
$ws.Cells.Item(4, 7).Value = 0.43
$ws.Cells.Item(5, 1).Value = "09-05-2025 01:00"
$ws.Cells.Item(5, 2).Value = "WORLD"
$ws.Cells.Item(5, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(5, 4).Value = "Velez Sarsfield - Olimpia"
$ws.Cells.Item(5, 5).Value = 50
$ws.Cells.Item(5, 6).Value = 1.73
$ws.Cells.Item(5, 7).Value = -0.14
$ws.Cells.Item(6, 1).Value = "08-05-2025 03:30"
$ws.Cells.Item(6, 2).Value = "WORLD"
$ws.Cells.Item(6, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(6, 4).Value = "Mushuc Runa SC - Cruzeiro"
$ws.Cells.Item(6, 5).Value = 70
$ws.Cells.Item(6, 6).Value = 2.2
$ws.Cells.Item(6, 7).Value = 0.54
$ws.Cells.Item(7, 1).Value = "08-05-2025 12:00"
$ws.Cells.Item(7, 2).Value = "CZECH-REPUBLIC"
$ws.Cells.Item(7, 3).Value = "4. LIGA - DIVIZIE C"
$ws.Cells.Item(7, 4).Value = "Benešov - Turnov"
$ws.Cells.Item(7, 5).Value = 60
$ws.Cells.Item(7, 6).Value = 1.73
$ws.Cells.Item(7, 7).Value = 0.04
$ws.Cells.Item(8, 1).Value = "08-05-2025 23:30"
$ws.Cells.Item(8, 2).Value = "ECUADOR"
$ws.Cells.Item(8, 3).Value = "LIGA PRO SERIE B"
$ws.Cells.Item(8, 4).Value = "Cumbayá - San Antonio"
$ws.Cells.Item(8, 5).Value = 50
$ws.Cells.Item(8, 6).Value = 2.05
$ws.Cells.Item(8, 7).Value = 0.02
$ws.Cells.Item(9, 1).Value = "09-05-2025 19:00"
$ws.Cells.Item(9, 2).Value = "CROATIA"
$ws.Cells.Item(9, 3).Value = "HNL"
$ws.Cells.Item(9, 4).Value = "NK Osijek - NK Lokomotiva Zagreb"
$ws.Cells.Item(9, 5).Value = 53.3
$ws.Cells.Item(9, 6).Value = 1.91
$ws.Cells.Item(9, 7).Value = 0.02
$ws.Cells.Item(10, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(10, 2).Value = "ITALY"
$ws.Cells.Item(10, 3).Value = "SERIE B"
$ws.Cells.Item(10, 4).Value = "Juve Stabia - Reggiana"
$ws.Cells.Item(10, 5).Value = 73.3
$ws.Cells.Item(10, 6).Value = 2.25
$ws.Cells.Item(10, 7).Value = 0.65
$ws.Cells.Item(11, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(11, 2).Value = "ITALY"
$ws.Cells.Item(11, 3).Value = "SERIE B"
$ws.Cells.Item(11, 4).Value = "Sassuolo - Catanzaro"
$ws.Cells.Item(11, 5).Value = 70
$ws.Cells.Item(11, 6).Value = 1.75
$ws.Cells.Item(11, 7).Value = 0.22
$ws.Cells.Item(12, 1).Value = "09-05-2025 21:00"
$ws.Cells.Item(12, 2).Value = "NETHERLANDS"
$ws.Cells.Item(12, 3).Value = "EERSTE DIVISIE"
$ws.Cells.Item(12, 4).Value = "Vitesse - Den Bosch"
$ws.Cells.Item(12, 5).Value = 60
$ws.Cells.Item(12, 6).Value = 2.35
$ws.Cells.Item(12, 7).Value = 0.41
$ws.Cells.Item(13, 1).Value = "09-05-2025 05:00"
$ws.Cells.Item(13, 2).Value = "WORLD"
$ws.Cells.Item(13, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(13, 4).Value = "Universitario - Independiente Del Valle"
$ws.Cells.Item(13, 5).Value = 66.7
$ws.Cells.Item(13, 6).Value = 2.15
$ws.Cells.Item(13, 7).Value = 0.43
$ws.Cells.Item(14, 1).Value = "09-05-2025 18:00"
$ws.Cells.Item(14, 2).Value = "ALGERIA"
$ws.Cells.Item(14, 3).Value = "LIGUE 2"
$ws.Cells.Item(14, 4).Value = "IB Khémis El Khechna - USM Annaba"
$ws.Cells.Item(14, 5).Value = 51.3
$ws.Cells.Item(14, 6).Value = 1.97
$ws.Cells.Item(14, 7).Value = 0.01
$ws.Cells.Item(15, 1).Value = "09-05-2025 18:00"
$ws.Cells.Item(15, 2).Value = "ALGERIA"
$ws.Cells.Item(15, 3).Value = "LIGUE 2"
$ws.Cells.Item(15, 4).Value = "Khroub - MB Rouisset"
$ws.Cells.Item(15, 5).Value = 50
$ws.Cells.Item(15, 6).Value = 2.45
$ws.Cells.Item(15, 7).Value = 0.23
$ws.Cells.Item(16, 1).Value = "09-05-2025 18:00"
$ws.Cells.Item(16, 2).Value = "ALGERIA"
$ws.Cells.Item(16, 3).Value = "LIGUE 2"
$ws.Cells.Item(16, 4).Value = "HB Chelghoum Laïd - Usm El Harrach"
$ws.Cells.Item(16, 5).Value = 60
$ws.Cells.Item(16, 6).Value = 2.63
$ws.Cells.Item(16, 7).Value = 0.58
$ws.Cells.Item(17, 1).Value = "09-05-2025 19:00"
$ws.Cells.Item(17, 2).Value = "AUSTRIA"
$ws.Cells.Item(17, 3).Value = "2. LIGA"
$ws.Cells.Item(17, 4).Value = "Floridsdorfer AC - SV Horn"
$ws.Cells.Item(17, 5).Value = 53.3
$ws.Cells.Item(17, 6).Value = 1.85
$ws.Cells.Item(17, 7).Value = -0.01
$ws.Cells.Item(18, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(18, 2).Value = "AUSTRIA"
$ws.Cells.Item(18, 3).Value = "2. LIGA"
$ws.Cells.Item(18, 4).Value = "Admira Wacker - Voitsberg"
$ws.Cells.Item(18, 5).Value = 50
$ws.Cells.Item(18, 6).Value = 1.77
$ws.Cells.Item(18, 7).Value = -0.12
$ws.Cells.Item(19, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(19, 2).Value = "AUSTRIA"
$ws.Cells.Item(19, 3).Value = "REGIONALLIGA - MITTE"
$ws.Cells.Item(19, 4).Value = "Treibach - Gleisdorf 09"
$ws.Cells.Item(19, 5).Value = 50
$ws.Cells.Item(19, 6).Value = 1.8
$ws.Cells.Item(19, 7).Value = -0.1
$ws.Cells.Item(20, 1).Value = "09-05-2025 21:45"
$ws.Cells.Item(20, 2).Value = "BELGIUM"
$ws.Cells.Item(20, 3).Value = "JUPILER PRO LEAGUE"
$ws.Cells.Item(20, 4).Value = "Charleroi - KVC Westerlo"
$ws.Cells.Item(20, 5).Value = 53.3
$ws.Cells.Item(20, 6).Value = 2
$ws.Cells.Item(20, 7).Value = 0.07000000000000001
$ws.Cells.Item(21, 1).Value = "09-05-2025 16:50"
$ws.Cells.Item(21, 2).Value = "CROATIA"
$ws.Cells.Item(21, 3).Value = "FIRST NL"
$ws.Cells.Item(21, 4).Value = "Sesvete - Croatia Zmijavci"
$ws.Cells.Item(21, 5).Value = 60
$ws.Cells.Item(21, 6).Value = 1.7
$ws.Cells.Item(21, 7).Value = 0.02
$ws.Cells.Item(22, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(22, 2).Value = "DENMARK"
$ws.Cells.Item(22, 3).Value = "2. DIVISION"
$ws.Cells.Item(22, 4).Value = "Fremad Amager - AB Copenhagen"
$ws.Cells.Item(22, 5).Value = 51.7
$ws.Cells.Item(22, 6).Value = 1.95
$ws.Cells.Item(22, 7).Value = 0.01
$ws.Cells.Item(23, 1).Value = "09-05-2025 20:30"
$ws.Cells.Item(23, 2).Value = "FRANCE"
$ws.Cells.Item(23, 3).Value = "NATIONAL 1"
$ws.Cells.Item(23, 4).Value = "Rouen - Gobelins"
$ws.Cells.Item(23, 5).Value = 73.3
$ws.Cells.Item(23, 6).Value = 2.45
$ws.Cells.Item(23, 7).Value = 0.8
$ws.Cells.Item(24, 1).Value = "09-05-2025 21:00"
$ws.Cells.Item(24, 2).Value = "HUNGARY"
$ws.Cells.Item(24, 3).Value = "NB I"
$ws.Cells.Item(24, 4).Value = "Gyori ETO FC - Fehérvár FC"
$ws.Cells.Item(24, 5).Value = 73.3
$ws.Cells.Item(24, 6).Value = 1.83
$ws.Cells.Item(24, 7).Value = 0.34
$ws.Cells.Item(25, 1).Value = "09-05-2025 11:30"
$ws.Cells.Item(25, 2).Value = "INDONESIA"
$ws.Cells.Item(25, 3).Value = "LIGA 1"
$ws.Cells.Item(25, 4).Value = "PSIS Semarang - PSS Sleman"
$ws.Cells.Item(25, 5).Value = 51.7
$ws.Cells.Item(25, 6).Value = 2.36
$ws.Cells.Item(25, 7).Value = 0.22
$ws.Cells.Item(26, 1).Value = "09-05-2025 21:45"
$ws.Cells.Item(26, 2).Value = "IRELAND"
$ws.Cells.Item(26, 3).Value = "FIRST DIVISION"
$ws.Cells.Item(26, 4).Value = "Finn Harps - Athlone Town"
$ws.Cells.Item(26, 5).Value = 50
$ws.Cells.Item(26, 6).Value = 2.5
$ws.Cells.Item(26, 7).Value = 0.25
$ws.Cells.Item(27, 1).Value = "10-05-2025 01:30"
$ws.Cells.Item(27, 2).Value = "PARAGUAY"
$ws.Cells.Item(27, 3).Value = "DIVISION PROFESIONAL - APERTURA"
$ws.Cells.Item(27, 4).Value = "Atlético Tembetary - Nacional Asuncion"
$ws.Cells.Item(27, 5).Value = 50
$ws.Cells.Item(27, 6).Value = 3.2
$ws.Cells.Item(27, 7).Value = 0.6
$ws.Cells.Item(28, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(28, 2).Value = "POLAND"
$ws.Cells.Item(28, 3).Value = "EKSTRAKLASA"
$ws.Cells.Item(28, 4).Value = "Gornik Zabrze - Slask Wroclaw"
$ws.Cells.Item(28, 5).Value = 53.3
$ws.Cells.Item(28, 6).Value = 2.05
$ws.Cells.Item(28, 7).Value = 0.09
$ws.Cells.Item(29, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(29, 2).Value = "POLAND"
$ws.Cells.Item(29, 3).Value = "I LIGA"
$ws.Cells.Item(29, 4).Value = "Tychy 71 - Wisla Krakow"
$ws.Cells.Item(29, 5).Value = 62.3
$ws.Cells.Item(29, 6).Value = 4
$ws.Cells.Item(29, 7).Value = 1.49
$ws.Cells.Item(30, 1).Value = "09-05-2025 19:00"
$ws.Cells.Item(30, 2).Value = "POLAND"
$ws.Cells.Item(30, 3).Value = "II LIGA - EAST"
$ws.Cells.Item(30, 4).Value = "Radunia Stężyca - Zaglebie Sosnowiec"
$ws.Cells.Item(30, 5).Value = 60
$ws.Cells.Item(30, 6).Value = 3
$ws.Cells.Item(30, 7).Value = 0.8
$ws.Cells.Item(31, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(31, 2).Value = "SPAIN"
$ws.Cells.Item(31, 3).Value = "SEGUNDA DIVISIÓN"
$ws.Cells.Item(31, 4).Value = "Cadiz - Almeria"
$ws.Cells.Item(31, 5).Value = 53.3
$ws.Cells.Item(31, 6).Value = 3
$ws.Cells.Item(31, 7).Value = 0.6
$ws.Cells.Item(32, 1).Value = "09-05-2025 20:30"
$ws.Cells.Item(32, 2).Value = "SWITZERLAND"
$ws.Cells.Item(32, 3).Value = "CHALLENGE LEAGUE"
$ws.Cells.Item(32, 4).Value = "Stade Lausanne-Ouchy - Bellinzona"
$ws.Cells.Item(32, 5).Value = 50
$ws.Cells.Item(32, 6).Value = 2.1
$ws.Cells.Item(32, 7).Value = 0.05
$ws.Cells.Item(33, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(33, 2).Value = "TURKEY"
$ws.Cells.Item(33, 3).Value = "SÜPER LIG"
$ws.Cells.Item(33, 4).Value = "Gazişehir Gaziantep - Alanyaspor"
$ws.Cells.Item(33, 5).Value = 63.3
$ws.Cells.Item(33, 6).Value = 2.38
$ws.Cells.Item(33, 7).Value = 0.51
$ws.Cells.Item(34, 1).Value = "09-05-2025 02:00"
$ws.Cells.Item(34, 2).Value = "USA"
$ws.Cells.Item(34, 3).Value = "USL CHAMPIONSHIP"
$ws.Cells.Item(34, 4).Value = "North Carolina - Orange County SC"
$ws.Cells.Item(34, 5).Value = 50
$ws.Cells.Item(34, 6).Value = 2.08
$ws.Cells.Item(34, 7).Value = 0.04
$ws.Cells.Item(35, 1).Value = "09-05-2025 21:00"
$ws.Cells.Item(35, 2).Value = "WORLD"
$ws.Cells.Item(35, 3).Value = "AFRICA CUP OF NATIONS U20"
$ws.Cells.Item(35, 4).Value = "South Africa U20 - Zambia U20"
$ws.Cells.Item(35, 5).Value = 80
$ws.Cells.Item(35, 6).Value = 2.3
$ws.Cells.Item(35, 7).Value = 0.84
$ws.Cells.Item(36, 1).Value = "09-05-2025 21:00"
$ws.Cells.Item(36, 2).Value = "WORLD"
$ws.Cells.Item(36, 3).Value = "AFRICA CUP OF NATIONS U20"
$ws.Cells.Item(36, 4).Value = "Tanzania U20 - Egypt U20"
$ws.Cells.Item(36, 5).Value = 60
$ws.Cells.Item(36, 6).Value = 7.5
$ws.Cells.Item(36, 7).Value = 3.5

# ---- Sheet: EV Away win ----
$ws = $wb.Worksheets.Item("EV Away win")
$ws.Cells.Item(2, 1).Value = "08-05-2025 03:30"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(2, 4).Value = "Cerro Porteno - Palmeiras"
$ws.Cells.Item(2, 5).Value = 66.7
$ws.Cells.Item(2, 6).Value = 1.85
$ws.Cells.Item(2, 7).Value = 0.23
$ws.Cells.Item(3, 1).Value = "09-05-2025 01:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(3, 4).Value = "Racing Montevideo - Huracan"
$ws.Cells.Item(3, 5).Value = 50
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(4, 1).Value = "08-05-2025 18:15"
$ws.Cells.Item(4, 2).Value = "OMAN"
$ws.Cells.Item(4, 3).Value = "PROFESSIONAL LEAGUE"
$ws.Cells.Item(4, 4).Value = "Bahla - Al-Shabab"
$ws.Cells.Item(4, 5).Value = 53.3
$ws.Cells.Item(4, 6).Value = 2.8
$ws.Cells.Item(4, 7).Value = 0.49
$ws.Cells.Item(5, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(5, 2).Value = "ITALY"
$ws.Cells.Item(5, 3).Value = "SERIE B"
$ws.Cells.Item(5, 4).Value = "Modena - Brescia"
$ws.Cells.Item(5, 5).Value = 60
$ws.Cells.Item(5, 6).Value = 3.2
$ws.Cells.Item(5, 7).Value = 0.92
$ws.Cells.Item(6, 1).Value = "09-05-2025 21:00"
$ws.Cells.Item(6, 2).Value = "NETHERLANDS"
$ws.Cells.Item(6, 3).Value = "EERSTE DIVISIE"
$ws.Cells.Item(6, 4).Value = "Jong AZ - Excelsior"
$ws.Cells.Item(6, 5).Value = 55.7
$ws.Cells.Item(6, 6).Value = 1.7
$ws.Cells.Item(6, 7).Value = -0.05
$ws.Cells.Item(7, 1).Value = "09-05-2025 05:00"
$ws.Cells.Item(7, 2).Value = "WORLD"
$ws.Cells.Item(7, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(7, 4).Value = "Cienciano - Caracas FC"
$ws.Cells.Item(7, 5).Value = 50
$ws.Cells.Item(7, 6).Value = 6.25
$ws.Cells.Item(7, 7).Value = 2.12
$ws.Cells.Item(8, 1).Value = "09-05-2025 19:00"
$ws.Cells.Item(8, 2).Value = "AUSTRIA"
$ws.Cells.Item(8, 3).Value = "2. LIGA"
$ws.Cells.Item(8, 4).Value = "Austria Lustenau - First Vienna"
$ws.Cells.Item(8, 5).Value = 50
$ws.Cells.Item(8, 6).Value = 2.4
$ws.Cells.Item(8, 7).Value = 0.2
$ws.Cells.Item(9, 1).Value = "09-05-2025 20:30"
$ws.Cells.Item(9, 2).Value = "AUSTRIA"
$ws.Cells.Item(9, 3).Value = "REGIONALLIGA - OST"
$ws.Cells.Item(9, 4).Value = "Austria Vienna (Am) - TWL Elektra"
$ws.Cells.Item(9, 5).Value = 66.7
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(9, 7).Value = 1.67
$ws.Cells.Item(10, 1).Value = "09-05-2025 19:00"
$ws.Cells.Item(10, 2).Value = "CYPRUS"
$ws.Cells.Item(10, 3).Value = "1. DIVISION"
$ws.Cells.Item(10, 4).Value = "Nea Salamis - Anorthosis"
$ws.Cells.Item(10, 5).Value = 60
$ws.Cells.Item(10, 6).Value = 2.05
$ws.Cells.Item(10, 7).Value = 0.23
$ws.Cells.Item(11, 1).Value = "09-05-2025 17:00"
$ws.Cells.Item(11, 2).Value = "EGYPT"
$ws.Cells.Item(11, 3).Value = "PREMIER LEAGUE"
$ws.Cells.Item(11, 4).Value = "Ceramica Cleopatra - Zamalek SC"
$ws.Cells.Item(11, 5).Value = 53.3
$ws.Cells.Item(11, 6).Value = 2.15
$ws.Cells.Item(11, 7).Value = 0.15
$ws.Cells.Item(12, 1).Value = "09-05-2025 19:30"
$ws.Cells.Item(12, 2).Value = "GERMANY"
$ws.Cells.Item(12, 3).Value = "2. BUNDESLIGA"
$ws.Cells.Item(12, 4).Value = "1. FC Nürnberg - 1.FC Köln"
$ws.Cells.Item(12, 5).Value = 50
$ws.Cells.Item(12, 6).Value = 1.77
$ws.Cells.Item(12, 7).Value = -0.12
$ws.Cells.Item(13, 1).Value = "09-05-2025 19:30"
$ws.Cells.Item(13, 2).Value = "GERMANY"
$ws.Cells.Item(13, 3).Value = "2. BUNDESLIGA"
$ws.Cells.Item(13, 4).Value = "Preußen Münster - Hertha BSC"
$ws.Cells.Item(13, 5).Value = 60
$ws.Cells.Item(13, 6).Value = 2.4
$ws.Cells.Item(13, 7).Value = 0.44
$ws.Cells.Item(14, 1).Value = "09-05-2025 19:30"
$ws.Cells.Item(14, 2).Value = "GERMANY"
$ws.Cells.Item(14, 3).Value = "REGIONALLIGA - BAYERN"
$ws.Cells.Item(14, 4).Value = "Vilzing - FC Schweinfurt 05"
$ws.Cells.Item(14, 5).Value = 60
$ws.Cells.Item(14, 6).Value = 2.6
$ws.Cells.Item(14, 7).Value = 0.5600000000000001
$ws.Cells.Item(15, 1).Value = "09-05-2025 20:30"
$ws.Cells.Item(15, 2).Value = "GERMANY"
$ws.Cells.Item(15, 3).Value = "REGIONALLIGA - NORD"
$ws.Cells.Item(15, 4).Value = "Hamburger SV II - Weiche Flensburg"
$ws.Cells.Item(15, 5).Value = 53.3
$ws.Cells.Item(15, 6).Value = 2.6
$ws.Cells.Item(15, 7).Value = 0.39
$ws.Cells.Item(16, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(16, 2).Value = "GERMANY"
$ws.Cells.Item(16, 3).Value = "REGIONALLIGA - NORDOST"
$ws.Cells.Item(16, 4).Value = "Hertha BSC II - BFC Dynamo"
$ws.Cells.Item(16, 5).Value = 60
$ws.Cells.Item(16, 6).Value = 2.63
$ws.Cells.Item(16, 7).Value = 0.58
$ws.Cells.Item(17, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(17, 2).Value = "LITHUANIA"
$ws.Cells.Item(17, 3).Value = "A LYGA"
$ws.Cells.Item(17, 4).Value = "Džiugas Telšiai - Kauno Žalgiris"
$ws.Cells.Item(17, 5).Value = 73.3
$ws.Cells.Item(17, 6).Value = 1.7
$ws.Cells.Item(17, 7).Value = 0.25
$ws.Cells.Item(18, 1).Value = "09-05-2025 23:30"
$ws.Cells.Item(18, 2).Value = "PERU"
$ws.Cells.Item(18, 3).Value = "PRIMERA DIVISIÓN"
$ws.Cells.Item(18, 4).Value = "Comerciantes Unidos - Deportivo Binacional"
$ws.Cells.Item(18, 5).Value = 53.3
$ws.Cells.Item(18, 6).Value = 3.6
$ws.Cells.Item(18, 7).Value = 0.92

# ---- Sheet: EV Over 2.5 ----
$ws = $wb.Worksheets.Item("EV Over 2.5")
$ws.Cells.Item(2, 1).Value = "09-05-2025 01:00"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(2, 4).Value = "Deportes Iquique - Atletico-MG"
$ws.Cells.Item(2, 5).Value = 80
$ws.Cells.Item(2, 6).Value = 1.95
$ws.Cells.Item(2, 7).Value = 0.5600000000000001
$ws.Cells.Item(3, 1).Value = "09-05-2025 05:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(3, 4).Value = "Universitario - Independiente Del Valle"
$ws.Cells.Item(3, 5).Value = 85
$ws.Cells.Item(3, 6).Value = 2.5
$ws.Cells.Item(3, 7).Value = 1.12
$ws.Cells.Item(4, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(4, 2).Value = "AUSTRIA"
$ws.Cells.Item(4, 3).Value = "2. LIGA"
$ws.Cells.Item(4, 4).Value = "Admira Wacker - Voitsberg"
$ws.Cells.Item(4, 5).Value = 80
$ws.Cells.Item(4, 6).Value = 1.76
$ws.Cells.Item(4, 7).Value = 0.41
$ws.Cells.Item(5, 1).Value = "09-05-2025 17:00"
$ws.Cells.Item(5, 2).Value = "EGYPT"
$ws.Cells.Item(5, 3).Value = "PREMIER LEAGUE"
$ws.Cells.Item(5, 4).Value = "Ceramica Cleopatra - Zamalek SC"
$ws.Cells.Item(5, 5).Value = 70
$ws.Cells.Item(5, 6).Value = 2.05
$ws.Cells.Item(5, 7).Value = 0.43
$ws.Cells.Item(6, 1).Value = "09-05-2025 21:45"
$ws.Cells.Item(6, 2).Value = "IRELAND"
$ws.Cells.Item(6, 3).Value = "PREMIER DIVISION"
$ws.Cells.Item(6, 4).Value = "Waterford - Galway United"
$ws.Cells.Item(6, 5).Value = 75
$ws.Cells.Item(6, 6).Value = 2.1
$ws.Cells.Item(6, 7).Value = 0.58
$ws.Cells.Item(7, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(7, 2).Value = "TURKEY"
$ws.Cells.Item(7, 3).Value = "SÜPER LIG"
$ws.Cells.Item(7, 4).Value = "Gazişehir Gaziantep - Alanyaspor"
$ws.Cells.Item(7, 5).Value = 78.8
$ws.Cells.Item(7, 6).Value = 1.7
$ws.Cells.Item(7, 7).Value = 0.34

# ---- Sheet: EV Btts ----
$ws = $wb.Worksheets.Item("EV Btts")
$ws.Cells.Item(2, 1).Value = "08-05-2025 03:30"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(2, 4).Value = "Atletico Grau - Gremio"
$ws.Cells.Item(2, 5).Value = 88
$ws.Cells.Item(2, 6).Value = 2.1
$ws.Cells.Item(2, 7).Value = 0.85
$ws.Cells.Item(3, 1).Value = "09-05-2025 01:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(3, 4).Value = "Deportes Iquique - Atletico-MG"
$ws.Cells.Item(3, 5).Value = 88
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 0.76
$ws.Cells.Item(4, 1).Value = "08-05-2025 20:00"
$ws.Cells.Item(4, 2).Value = "DENMARK"
$ws.Cells.Item(4, 3).Value = "DBU POKALEN"
$ws.Cells.Item(4, 4).Value = "FC Copenhagen - Viborg"
$ws.Cells.Item(4, 5).Value = 76.7
$ws.Cells.Item(4, 6).Value = 1.8
$ws.Cells.Item(4, 7).Value = 0.38
$ws.Cells.Item(5, 1).Value = "08-05-2025 18:00"
$ws.Cells.Item(5, 2).Value = "GEORGIA"
$ws.Cells.Item(5, 3).Value = "EROVNULI LIGA"
$ws.Cells.Item(5, 4).Value = "Dinamo Tbilisi - Telavi"
$ws.Cells.Item(5, 5).Value = 73.3
$ws.Cells.Item(5, 6).Value = 2.3
$ws.Cells.Item(5, 7).Value = 0.6899999999999999
$ws.Cells.Item(6, 1).Value = "08-05-2025 19:00"
$ws.Cells.Item(6, 2).Value = "NORWAY"
$ws.Cells.Item(6, 3).Value = "NM CUPEN"
$ws.Cells.Item(6, 4).Value = "Tromso - KFUM Oslo"
$ws.Cells.Item(6, 5).Value = 70
$ws.Cells.Item(6, 6).Value = 1.83
$ws.Cells.Item(6, 7).Value = 0.28
$ws.Cells.Item(7, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(7, 2).Value = "ITALY"
$ws.Cells.Item(7, 3).Value = "SERIE B"
$ws.Cells.Item(7, 4).Value = "Modena - Brescia"
$ws.Cells.Item(7, 5).Value = 83.3
$ws.Cells.Item(7, 6).Value = 1.91
$ws.Cells.Item(7, 7).Value = 0.59
$ws.Cells.Item(8, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(8, 2).Value = "ITALY"
$ws.Cells.Item(8, 3).Value = "SERIE B"
$ws.Cells.Item(8, 4).Value = "Palermo - Frosinone"
$ws.Cells.Item(8, 5).Value = 69.2
$ws.Cells.Item(8, 6).Value = 1.77
$ws.Cells.Item(8, 7).Value = 0.22
$ws.Cells.Item(9, 1).Value = "09-05-2025 21:30"
$ws.Cells.Item(9, 2).Value = "ITALY"
$ws.Cells.Item(9, 3).Value = "SERIE B"
$ws.Cells.Item(9, 4).Value = "Pisa - Sudtirol"
$ws.Cells.Item(9, 5).Value = 65
$ws.Cells.Item(9, 6).Value = 1.95
$ws.Cells.Item(9, 7).Value = 0.27
$ws.Cells.Item(10, 1).Value = "09-05-2025 05:00"
$ws.Cells.Item(10, 2).Value = "WORLD"
$ws.Cells.Item(10, 3).Value = "CONMEBOL LIBERTADORES"
$ws.Cells.Item(10, 4).Value = "Universitario - Independiente Del Valle"
$ws.Cells.Item(10, 5).Value = 73.3
$ws.Cells.Item(10, 6).Value = 2.1
$ws.Cells.Item(10, 7).Value = 0.54
$ws.Cells.Item(11, 1).Value = "09-05-2025 05:00"
$ws.Cells.Item(11, 2).Value = "WORLD"
$ws.Cells.Item(11, 3).Value = "CONMEBOL SUDAMERICANA"
$ws.Cells.Item(11, 4).Value = "Cienciano - Caracas FC"
$ws.Cells.Item(11, 5).Value = 68
$ws.Cells.Item(11, 6).Value = 1.91
$ws.Cells.Item(11, 7).Value = 0.3
$ws.Cells.Item(12, 1).Value = "09-05-2025 18:00"
$ws.Cells.Item(12, 2).Value = "ALGERIA"
$ws.Cells.Item(12, 3).Value = "LIGUE 2"
$ws.Cells.Item(12, 4).Value = "MSP Batna - US Souf"
$ws.Cells.Item(12, 5).Value = 72
$ws.Cells.Item(12, 6).Value = 2.25
$ws.Cells.Item(12, 7).Value = 0.62
$ws.Cells.Item(13, 1).Value = "09-05-2025 18:00"
$ws.Cells.Item(13, 2).Value = "ALGERIA"
$ws.Cells.Item(13, 3).Value = "LIGUE 2"
$ws.Cells.Item(13, 4).Value = "IB Khémis El Khechna - USM Annaba"
$ws.Cells.Item(13, 5).Value = 71.09999999999999
$ws.Cells.Item(13, 6).Value = 1.8
$ws.Cells.Item(13, 7).Value = 0.28
$ws.Cells.Item(14, 1).Value = "09-05-2025 12:30"
$ws.Cells.Item(14, 2).Value = "AUSTRALIA"
$ws.Cells.Item(14, 3).Value = "VICTORIA NPL"
$ws.Cells.Item(14, 4).Value = "Heidelberg United - Port Melbourne"
$ws.Cells.Item(14, 5).Value = 83.3
$ws.Cells.Item(14, 6).Value = 1.7
$ws.Cells.Item(14, 7).Value = 0.42
$ws.Cells.Item(15, 1).Value = "09-05-2025 20:00"
$ws.Cells.Item(15, 2).Value = "EGYPT"
$ws.Cells.Item(15, 3).Value = "PREMIER LEAGUE"
$ws.Cells.Item(15, 4).Value = "National Bank Of Egypt - Pyramids FC"
$ws.Cells.Item(15, 5).Value = 73.3
$ws.Cells.Item(15, 6).Value = 1.86
$ws.Cells.Item(15, 7).Value = 0.36
$ws.Cells.Item(16, 1).Value = "09-05-2025 21:45"
$ws.Cells.Item(16, 2).Value = "IRELAND"
$ws.Cells.Item(16, 3).Value = "FIRST DIVISION"
$ws.Cells.Item(16, 4).Value = "UCD - Longford Town"
$ws.Cells.Item(16, 5).Value = 66.7
$ws.Cells.Item(16, 6).Value = 1.9
$ws.Cells.Item(16, 7).Value = 0.27
$ws.Cells.Item(17, 1).Value = "09-05-2025 21:45"
$ws.Cells.Item(17, 2).Value = "IRELAND"
$ws.Cells.Item(17, 3).Value = "PREMIER DIVISION"
$ws.Cells.Item(17, 4).Value = "Cork City - Derry City"
$ws.Cells.Item(17, 5).Value = 70
$ws.Cells.Item(17, 6).Value = 1.95
$ws.Cells.Item(17, 7).Value = 0.36
$ws.Cells.Item(18, 1).Value = "09-05-2025 21:45"
$ws.Cells.Item(18, 2).Value = "IRELAND"
$ws.Cells.Item(18, 3).Value = "PREMIER DIVISION"
$ws.Cells.Item(18, 4).Value = "Waterford - Galway United"
$ws.Cells.Item(18, 5).Value = 73.3
$ws.Cells.Item(18, 6).Value = 1.91
$ws.Cells.Item(18, 7).Value = 0.4
$ws.Cells.Item(19, 1).Value = "10-05-2025 01:30"
$ws.Cells.Item(19, 2).Value = "PARAGUAY"
$ws.Cells.Item(19, 3).Value = "DIVISION PROFESIONAL - APERTURA"
$ws.Cells.Item(19, 4).Value = "Atlético Tembetary - Nacional Asuncion"
$ws.Cells.Item(19, 5).Value = 68
$ws.Cells.Item(19, 6).Value = 2.25
$ws.Cells.Item(19, 7).Value = 0.53
$ws.Cells.Item(20, 1).Value = "09-05-2025 23:00"
$ws.Cells.Item(20, 2).Value = "PARAGUAY"
$ws.Cells.Item(20, 3).Value = "DIVISION PROFESIONAL - APERTURA"
$ws.Cells.Item(20, 4).Value = "Sportivo Trinidense - 2 De Mayo"
$ws.Cells.Item(20, 5).Value = 73.3
$ws.Cells.Item(20, 6).Value = 2.2
$ws.Cells.Item(20, 7).Value = 0.61
$ws.Cells.Item(21, 1).Value = "09-05-2025 19:00"
$ws.Cells.Item(21, 2).Value = "POLAND"
$ws.Cells.Item(21, 3).Value = "EKSTRAKLASA"
$ws.Cells.Item(21, 4).Value = "Motor Lublin - Piast Gliwice"
$ws.Cells.Item(21, 5).Value = 80
$ws.Cells.Item(21, 6).Value = 1.77
$ws.Cells.Item(21, 7).Value = 0.42
$ws.Cells.Item(22, 1).Value = "09-05-2025 02:00"
$ws.Cells.Item(22, 2).Value = "USA"
$ws.Cells.Item(22, 3).Value = "USL CHAMPIONSHIP"
$ws.Cells.Item(22, 4).Value = "North Carolina - Orange County SC"
$ws.Cells.Item(22, 5).Value = 84
$ws.Cells.Item(22, 6).Value = 1.74
$ws.Cells.Item(22, 7).Value = 0.46

